$wb = $excel.ActiveWorkbook

$oldGuid = "507765fa-b26b-4afc-8831-365945693fc2"
$newGuid = "0450d315-3fd3-4ec8-92d1-515e80dfc30f"
$oldHash = "fe7e2fa703dab2db7ad16af9026bc89c204d696b"
$newHash = "9f240ad1406db5e4135481e7d4d18cb9fb6f5ebb"

# Original "HyperLink" font color (6495ED) expressed as the VBA-style
# R + G*256 + B*65536 integer that Font.Color expects.
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("D2").Value = "2016-50-18 20:50:16"

$urlOverviewA2 = "https://github.com/OpenLocalizationTest/oltest/blob/7e17630268bffdbcc5668f9b98f2294ce3c5b622/e2e/$oldGuid.md"

# The headless hyperlink model only supports appending a new hyperlink or
# wiping every hyperlink on the sheet, so clear everything and recreate each
# hyperlink with its original target but the refreshed display text.
$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $urlOverviewA2, "", "", "$newGuid.md")

# Re-apply the original hyperlink look (Hyperlinks.Add resets it to the
# built-in "Hyperlink" style/theme color otherwise).
$wsOverview.Range("A2").Font.Underline = $true
$wsOverview.Range("A2").Font.Color = $hyperlinkColor
$wsOverview.Range("A2").Font.Name = "Calibri"
$wsOverview.Range("A2").Font.Size = 11

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("E2").Value = "2016-03-18 20:50:13"

$urlZhCnA2 = "https://github.com/OpenLocalizationTest/oltest/blob/7e17630268bffdbcc5668f9b98f2294ce3c5b622/e2e/$oldGuid.md"
$urlZhCnB2 = "https://github.com/OpenLocalizationTest/oltest/blob/7e17630268bffdbcc5668f9b98f2294ce3c5b622/e2e/$oldGuid.md"
$urlZhCnD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a569000791b6028535f30feb2b24a4729b925666/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $urlZhCnA2, "", "", "$newGuid.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $urlZhCnB2, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $urlZhCnD2, "", "", "$newGuid.$newHash.zh-cn.xlf")

foreach ($addr in @("A2", "B2", "D2")) {
    $wsZhCn.Range($addr).Font.Underline = $true
    $wsZhCn.Range($addr).Font.Color = $hyperlinkColor
    $wsZhCn.Range($addr).Font.Name = "Calibri"
    $wsZhCn.Range($addr).Font.Size = 11
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("E2").Value = "2016-03-18 20:50:16"

$urlDeDeA2 = "https://github.com/OpenLocalizationTest/oltest/blob/7e17630268bffdbcc5668f9b98f2294ce3c5b622/e2e/$oldGuid.md"
$urlDeDeB2 = "https://github.com/OpenLocalizationTest/oltest/blob/7e17630268bffdbcc5668f9b98f2294ce3c5b622/e2e/$oldGuid.md"
$urlDeDeD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5fc78435dcbfb825e70619befff7d7669647f31e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $urlDeDeA2, "", "", "$newGuid.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $urlDeDeB2, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $urlDeDeD2, "", "", "$newGuid.$newHash.de-de.xlf")

foreach ($addr in @("A2", "B2", "D2")) {
    $wsDeDe.Range($addr).Font.Underline = $true
    $wsDeDe.Range($addr).Font.Color = $hyperlinkColor
    $wsDeDe.Range($addr).Font.Name = "Calibri"
    $wsDeDe.Range($addr).Font.Size = 11
}
